# Update the cryptocurrency price/volume snapshot data (row 2-51)
# to the latest values pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.354.15"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'1.844.83"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("D4").Value = "'0.9979"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = "'240.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = "'0.6266"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = "'0.9986"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = "'0.07488"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.82%  '
$ws.Range("D9").Value = "'0.2900"
$ws.Range("D9").ClearFormats()
$ws.Range("E10").Value = '  -1.42%  '
$ws.Range("D11").Value = "'0.07713"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").Value = "'1.844.35"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.32%  '
$ws.Range("D13").Value = "'5.000"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").Value = "'0.6786"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").Value = "'0.00001029"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.96%  '
$ws.Range("D16").Value = "'82.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("D17").Value = "'2.099.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.99%  '
$ws.Range("D18").Value = "'6.163"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").Value = "'29.398.43"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = "'229.14"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.00%  '
$ws.Range("D21").Value = "'12.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = "'0.9982"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = "'7.476"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = "'0.9995"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = "'158.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").Value = "'8.407"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = "'17.54"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("D29").Value = "'0.06496"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +16.17%  '
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("D31").Value = "'1.469"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").Value = "'4.058"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").Value = "'1.823"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("D36").Value = "'0.6991"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.55%  '
$ws.Range("D37").Value = "'2.575"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").Value = "'1.258.45"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.36%  '
$ws.Range("D39").Value = "'2.829"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.08%  '
$ws.Range("D40").Value = "'0.01827"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("D41").Value = "'6.588"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.13%  '
$ws.Range("D42").Value = "'0.9124"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").Value = "'0.9982"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("D44").Value = "'2.007.08"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -18.44%  '
$ws.Range("D45").Value = "'101.46"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").Value = "'66.10"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = "'1.735"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.31%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = "'7.073"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = "'0.1175"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = "'0.00000000116"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.09%  '
$ws.Range("D51").Value = "'9.003"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.15%  '

Write-Host "Updated cryptos list"
